# "bring back input tables in development"
# Re-introduce the scenario rows (id_scenario 10/20/30) that were missing from
# Sheet1's Table2, and add a "note" worksheet documenting the assumptions
# behind each scenario id.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------------------
# 1. Add the new rows of scenario data to Sheet1 (rows 4-9)
# ---------------------------------------------------------------------------

# Columns E..AT (price curve, 2010..2051) for the three distinct price
# trajectories used by scenarios 10, 20 and 30 (scenario 10 and 20 each use
# the same curve for both id_sector 3 and 6; likewise for 30).
$curveA = @(0,0,0,0,0,0,0,0,0,0,0,25,30,35,45,55,57.095061376267388,68.86542787537482,80.142141689967445,90.940240780826201,101.27435905137465,111.15873649478881,120.60722909684587,129.63331850022271,138.25012143582575,146.47039892660456,154.3065652691779,161.77069679847767,168.87454044049778,175.62952205811769,182.04675459485659,188.13704602130269,193.91090708885301,199.37855889529283,204.54994026663948,209.43471495957414,214.04227868868477,218.38176598264715,222.46205687337556,226.29178342208175,229.87933608609032,230)
$curveB = @(0,0,0,0,0,0,0,0,0,0,0,25,30,35,45,55,57.047530688133691,64.43271393768741,71.071070844983723,77.470120390413101,83.137179525687316,89.079368247394399,95.803614548422928,102.31665925011136,109.12506071791287,115.23519946330228,119.05745669816191,123.96437822031805,128.65413552053596,133.13557261086527,137.41696322588535,141.50608724573311,145.41029152080606,149.13653896626627,152.6914487051998,156.08132928927245,159.31220650001273,163.2796322619476,166.19244123194196,168.96374955227228,172.43966804304517,172.5)
$curveC = @(0,0,0,0,0,0,0,0,0,0,0,25,30,35,45,55,57,60,62,64,65,67,71,75,80,84,83.808348127145933,86.158059642158435,88.433730600574123,90.641623163612849,92.787171856914142,94.875128470163503,96.909675952759088,98.894519037239718,100.83295714376011,102.72794361897078,104.58213431134067,108.17749854124806,109.92282559050837,111.63571568246282,115,115)

$newRows = @(
    @{ Row = 4; A = 10; B = 9; C = 3; Curve = $curveA },
    @{ Row = 5; A = 10; B = 9; C = 6; Curve = $curveA },
    @{ Row = 6; A = 20; B = 9; C = 3; Curve = $curveB },
    @{ Row = 7; A = 20; B = 9; C = 6; Curve = $curveB },
    @{ Row = 8; A = 30; B = 9; C = 3; Curve = $curveC },
    @{ Row = 9; A = 30; B = 9; C = 6; Curve = $curveC }
)

foreach ($rowInfo in $newRows) {
    $r = $rowInfo.Row
    $sheet1.Cells.Item($r, 1).Value = $rowInfo.A
    $sheet1.Cells.Item($r, 2).Value = $rowInfo.B
    $sheet1.Cells.Item($r, 3).Value = $rowInfo.C
    $sheet1.Cells.Item($r, 4).Value = "euro/tCO2"

    $curve = $rowInfo.Curve
    for ($i = 0; $i -lt $curve.Length; $i++) {
        $cell = $sheet1.Cells.Item($r, 5 + $i)
        $cell.Value = $curve[$i]
        $cell.NumberFormat = "0"
    }
}

# ---------------------------------------------------------------------------
# 2. Resize Table2 (the ListObject) so it covers the newly added rows
# ---------------------------------------------------------------------------
$table = $sheet1.ListObjects.Item(1)
$table.Resize($sheet1.Range("A1:AT9"))

# Restore the view/selection on Sheet1
$sheet1.Activate()
[void]$sheet1.Range("V13").Select()

# ---------------------------------------------------------------------------
# 3. Add the "note" worksheet right after Sheet1, documenting assumptions
# ---------------------------------------------------------------------------
$note = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $sheet1)
$note.Name = "note"

$note.Cells.Item(1, 1).Value = "note: scenarios differentiate after 2026."

$note.Cells.Item(2, 1).Value = "id_scenario"
$note.Cells.Item(2, 2).Value = "assumption"

# Fill in the assumption column in the same order the original author did,
# so that newly created shared strings land at the expected indices.
$note.Cells.Item(5, 1).Value = 20
$note.Cells.Item(5, 2).Value = "average of 10 & 30 after 2026"

$note.Cells.Item(6, 1).Value = 30
$note.Cells.Item(6, 2).Value = "same as 1"

$note.Cells.Item(4, 1).Value = 10
$note.Cells.Item(4, 2).Value = 'same as "BEHG" (https://www.axpo.com/de/en/business/gas/behg.html)'

$note.Cells.Item(3, 1).Value = 1
$note.Cells.Item(3, 2).Value = "a slight steady increase after the announced price corridor (until end of 2026)"

[void]$note.Range("B4").Select()

# Keep Sheet1 as the active/visible tab, matching the target workbook state.
$sheet1.Activate()
